$wb = $excel.ActiveWorkbook

# --- Sheet "Przerobione": append 5 new vocabulary rows ---
$przerobione = $wb.Worksheets.Item("Przerobione")

$newWords = @(
    @(1042, "音乐", "yīnyuè", "Music"),
    @(405,  "继续", "jìxù",   "Continue"),
    @(1078, "元",   "yuán",   "element, Yuan"),
    @(828,  "所以", "suǒyǐ",  "therefore"),
    @(200,  "短",   "duǎn",   "short")
)

$startRow = 57
for ($i = 0; $i -lt $newWords.Count; $i++) {
    $row = $startRow + $i
    $word = $newWords[$i]
    $przerobione.Cells.Item($row, 1).Value = $word[0]
    $przerobione.Cells.Item($row, 2).Value = $word[1]
    $przerobione.Cells.Item($row, 3).Value = $word[2]
    $przerobione.Cells.Item($row, 4).Value = $word[3]
}

# --- Sheet "5 losowych": replace the 5 random words with the new ones ---
$losowych = $wb.Worksheets.Item("5 losowych")

for ($i = 0; $i -lt $newWords.Count; $i++) {
    $row = 2 + $i
    $word = $newWords[$i]
    $losowych.Cells.Item($row, 1).Value = $word[0]
    $losowych.Cells.Item($row, 2).Value = $word[1]
    $losowych.Cells.Item($row, 3).Value = $word[2]
    $losowych.Cells.Item($row, 4).Value = $word[3]
}
